$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: apply "hyperlink look" (blue, underlined, size 11, Calibri) to the
# full text of a cell while forcing Excel to keep the run as a *rich-text*
# shared-string run (a single Characters() call spanning 100% of the string
# collapses into a cell-level font style instead, so we always split into two
# adjoining runs that together cover the whole string).
# ---------------------------------------------------------------------------
function Format-LinkText($cell, [string]$text) {
    $n = $text.Length
    if ($n -lt 2) {
        $c1 = $cell.Characters(1, $n)
        $c1.Font.Underline = 2
        $c1.Font.Size = 11
        $c1.Font.ColorIndex = 11
        $c1.Font.Name = "Calibri"
        return
    }
    $half = [Math]::Floor($n / 2)
    $c1 = $cell.Characters(1, $half)
    $c1.Font.Underline = 2
    $c1.Font.Size = 11
    $c1.Font.ColorIndex = 11
    $c1.Font.Name = "Calibri"
    $c2 = $cell.Characters($half + 1, $n - $half)
    $c2.Font.Underline = 2
    $c2.Font.Size = 11
    $c2.Font.ColorIndex = 11
    $c2.Font.Name = "Calibri"
}

# ---------------------------------------------------------------------------
# Helper: fill one event row (columns A-E) using row 627 (a fully populated,
# already-correctly-styled row) as the formatting donor, then set the real
# values/hyperlink on top.
# ---------------------------------------------------------------------------
function Add-EventRow {
    param(
        [int]$RowNum,
        [int]$DateSerial,
        [string]$EventName,
        [string]$LocationName,
        [string]$CityName,
        [string]$LinkUrl,
        [bool]$BStyle10 = $false
    )

    # Clone formatting (styles incl. borders/fill/number-format) from the
    # template row so every cell lands on style "3" (text) / "4" (date).
    $ws.Range("A627:E627").Copy($ws.Range("A" + $RowNum + ":E" + $RowNum))

    $ws.Range("A$RowNum").Value = $DateSerial
    $ws.Range("B$RowNum").Value = $EventName
    $ws.Range("C$RowNum").Value = $LocationName
    $ws.Range("D$RowNum").Value = $CityName

    if ($BStyle10) {
        # A handful of source rows use a distinct font style ("10") on the
        # Event column; replicate it from a known donor (B555) then restore
        # the text (copy only brings the format across).
        $ws.Range("B555").Copy()
        $ws.Range("B$RowNum").PasteSpecial(-4122)
        $ws.Range("B$RowNum").Value = $EventName
    }

    $linkCell = $ws.Range("E$RowNum")
    $linkCell.Value = $LinkUrl
    Format-LinkText $linkCell $LinkUrl
    $ws.Hyperlinks.Add($linkCell, $LinkUrl, "", "", $LinkUrl) | Out-Null

    # Hyperlinks.Add swaps the cell onto a synthetic "hyperlink" style; put
    # it back onto style "3" (matching every other populated Link cell) by
    # re-pasting just the formatting from the template row's E cell.
    $ws.Range("E627").Copy()
    $linkCell.PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# 1) Row 608: date correction (45865 -> 45878)
# ---------------------------------------------------------------------------
$ws.Range("A608").Value = 45878

# ---------------------------------------------------------------------------
# 2) Rows 628-642: newly added events (previously blank placeholder rows)
# ---------------------------------------------------------------------------
Add-EventRow -RowNum 628 -DateSerial 45864 -EventName "SODOM UND GOMORRA x EHRENLOS" -LocationName "Triptychon" -CityName "Münster" -LinkUrl "https://www.instagram.com/reel/DMPW8TJowBQ/?igsh=cXF2MmhsdXZ0emNq" -BStyle10 $false
Add-EventRow -RowNum 629 -DateSerial 45862 -EventName "SALON ELECTRONIQUE" -LocationName "SNRS" -CityName "Dortmund" -LinkUrl "https://www.instagram.com/p/DMQUuqKIHQ-/?igsh=amxzaTA5ZzV1NWUz" -BStyle10 $false
Add-EventRow -RowNum 630 -DateSerial 45863 -EventName "CTRL" -LocationName "SNRS" -CityName "Dortmund" -LinkUrl "https://www.instagram.com/reel/DLsBp86M5Xk/?igsh=ZW92ZHlwaGdpY3Zy" -BStyle10 $false
Add-EventRow -RowNum 631 -DateSerial 45948 -EventName "HADES x PUMP" -LocationName "SNRS" -CityName "Dortmund" -LinkUrl "https://www.instagram.com/reel/DMalZU5sfFM/?igsh=MWtyMnEycmRxdW9pcA==" -BStyle10 $false
Add-EventRow -RowNum 632 -DateSerial 45884 -EventName "VROLIK B-DAY BASH" -LocationName "Zimmermanns" -CityName "Köln" -LinkUrl "https://www.instagram.com/reel/DLMyoRtoLUf/?igsh=YzZkMmtqeHUyYWRw" -BStyle10 $true
Add-EventRow -RowNum 633 -DateSerial 45864 -EventName "SOMMERFEST DAY & NIGHT" -LocationName "Rotunde" -CityName "Bochum" -LinkUrl "https://www.instagram.com/reel/DMTQnz6sqzi/?igsh=MXZmbDhodTY5ZGpxNw==" -BStyle10 $false
Add-EventRow -RowNum 634 -DateSerial 45864 -EventName "RAVE IM REINEKE" -LocationName "Reineke Fuchs" -CityName "Köln" -LinkUrl "https://www.instagram.com/reel/DMC_-UVMvbl/?igsh=MW0xcHR6YWw4bnZpcQ==" -BStyle10 $false
Add-EventRow -RowNum 635 -DateSerial 45885 -EventName "EISKALT TECHNO BRUNCH" -LocationName "check event link" -CityName "Köln" -LinkUrl "https://www.instagram.com/reel/DMa0sA3CyRL/?igsh=MXR2OHl1bXNkY2E5eg==" -BStyle10 $false
Add-EventRow -RowNum 636 -DateSerial 45864 -EventName "BZZ BZZ TECHNO" -LocationName "Helios37" -CityName "Köln" -LinkUrl "https://www.instagram.com/reel/DMaeCLyImLK/?igsh=dTM3YXljb21uMDdn" -BStyle10 $false
Add-EventRow -RowNum 637 -DateSerial 45871 -EventName "BZZ BZZ TECHNO" -LocationName "Zimmermanns" -CityName "Köln" -LinkUrl "https://www.instagram.com/reel/DMP4xS4s6uH/?igsh=aWpxMjh0cGtiNWhk" -BStyle10 $false
Add-EventRow -RowNum 638 -DateSerial 45861 -EventName "#MITTWOCHENENDE" -LocationName "Odonien" -CityName "Köln" -LinkUrl "https://www.instagram.com/p/DMN7rWUsZlF/?igsh=MXNndzJseDUwamRleQ==" -BStyle10 $false
Add-EventRow -RowNum 639 -DateSerial 45877 -EventName "PROJEKT RAVE" -LocationName "Sam‘s" -CityName "Bielefeld" -LinkUrl "https://www.instagram.com/reel/DMIuqDOsB9Z/?igsh=MWk3OHoxZXhxemEyZA==" -BStyle10 $false
Add-EventRow -RowNum 640 -DateSerial 45870 -EventName "NOCTURGENERATION & ARTFACTORY" -LocationName "check event link" -CityName "Essen" -LinkUrl "https://www.instagram.com/reel/DMSG9eRN4b_/?igsh=MXM0bGdyeWJzanVkZg==" -BStyle10 $false
Add-EventRow -RowNum 641 -DateSerial 45864 -EventName "HAPPY HIPPIE OPEN AIR" -LocationName "Brauerei" -CityName "Bielefeld" -LinkUrl "https://www.instagram.com/p/DL4KZc8tfcV/?igsh=MXU5dG1pYzNocmM0bg==" -BStyle10 $false
Add-EventRow -RowNum 642 -DateSerial 45871 -EventName "WYLDHEARTS OPEN AIR" -LocationName "Fühlinger See" -CityName "Köln" -LinkUrl "https://www.instagram.com/reel/DMDxyoHt0Gx/?igsh=M2ZrZ3kzazE2Znpo" -BStyle10 $false

# ---------------------------------------------------------------------------
# 3) Rows 643-655: extend the trailing blank placeholder block (same style
#    as every other empty row: date col style "4", rest style "5") so the
#    sheet's used range grows to A1:E655.
# ---------------------------------------------------------------------------
for ($r = 643; $r -le 655; $r++) {
    $ws.Range("A639:E639").Copy($ws.Range("A" + $r + ":E" + $r))
    $ws.Rows.Item($r).RowHeight = 15
}
